$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy cell formatting from column R into column S for rows 3-37 (row 34 has no value,
# it only gains a blank formatted cell S34 to mirror P34/Q34/R34).
$ws.Range("R3:R34").Copy() | Out-Null
$ws.Range("S3:S34").PasteSpecial(-4122) | Out-Null
$ws.Range("R35:R37").Copy() | Out-Null
$ws.Range("S35:S37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new 2022 column (S) values.
# (values are written in plain decimal form -- the interpreter's expression
# parser does not accept scientific-notation numeric literals like 1E-2)
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 0.15686557910355481
$ws.Range("S5").Value = 0.18747863920572591
$ws.Range("S6").Value = 0.12556588018347117
$ws.Range("S7").Value = 0.051313356512815066
$ws.Range("S8").Value = 0.060745743331368028
$ws.Range("S9").Value = 0.042060988433228183
$ws.Range("S10").Value = 0.092022006630303563
$ws.Range("S11").Value = 0.078942235953699605
$ws.Range("S12").Value = 0.10098382728705417
$ws.Range("S13").Value = 0.097010038673425045
$ws.Range("S14").Value = 0.12657756598786343
$ws.Range("S15").Value = 0.067310604785784003
$ws.Range("S16").Value = 0.12618253497302423
$ws.Range("S17").Value = 0.15767275020694549
$ws.Range("S18").Value = 0.083781780685077176
$ws.Range("S19").Value = 0.089790167285988584
$ws.Range("S20").Value = 0.11543537913568107
$ws.Range("S21").Value = 0.064489306438090949
$ws.Range("S22").Value = 0.077235413540471365
$ws.Range("S23").Value = 0.1335826876836762
$ws.Range("S24").Value = 0.021874179718260566
$ws.Range("S25").Value = 0.13849188927432132
$ws.Range("S26").Value = 0.15541703258327452
$ws.Range("S27").Value = 0.12135301021830269
$ws.Range("S28").Value = 0.4304881257025327
$ws.Range("S29").Value = 0.49554896622979544
$ws.Range("S30").Value = 0.35193780867878632
$ws.Range("S31").Value = 0.21076296192215821
$ws.Range("S32").Value = 0.25905990040586052
$ws.Range("S33").Value = 0.1647039446594746
$ws.Range("S35").Value = 0
$ws.Range("S36").Value = 0.1
$ws.Range("S37").Value = 0.2

# Update the active cell / selection to match the saved view state.
$ws.Range("T15").Select() | Out-Null
